# Negate every numeric value in column E ("Block") across the data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    $val = $cell.Value2
    if ($val -is [double]) {
        $cell.Value = -$val
    }
}
